# Leave Card update - 6/22/2023 5:35 PM
# Inserts a new SL(3-0-0) leave entry row into the leave ledger table
# (Table1) at row 78, shifting the existing period rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# Insert a brand-new row into the worksheet at row 78, shifting every
# row below (the remaining ledger rows + the table's closing row) down
# by one.
$ws.Rows.Item(78).Insert(-4121)   # xlShiftDown

# Grow the table definition so it covers the newly inserted row at the
# bottom (the physical rows already shifted down, so the table now
# needs to include one more row).
$lo.Resize($ws.Range("A8:K133"))

# The freshly inserted row 78 has no formatting yet; pick up the same
# look as the rest of the ledger rows (copy from row 79, which holds
# what used to be row 78 before the insert).
$ws.Range("A79:K79").Copy()
$ws.Range("A78:K78").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Restore/ensure the calculated column formulas on the new row and on
# the row that is now the last table row (the insert operation can
# leave these blank / stale).
$ws.Range("G78").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G133").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Fill in the new leave record values.
$ws.Range("B78").Value = "SL(3-0-0)"
$ws.Range("H78").Value = 3
$ws.Range("K78").Value = "6/9,13,15/2023"

# Update the active selection to reflect where the user ended up
# working (K79, the cell that used to be the selection anchor K78
# before the insert shifted it down).
$ws.Range("K79").Select()
